$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 and J1 ("I0" and "IF"), matching the formatting
# of the existing header row (style index 1 -> bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data for columns I and J, rows 2 through 58
$rowsData = @(
    @{Row=2; I=9; J=9}
    @{Row=3; I=9; J=9}
    @{Row=4; I=8; J=8}
    @{Row=5; I=7; J=7}
    @{Row=6; I=8; J=8}
    @{Row=7; I=8; J=8}
    @{Row=8; I=8; J=8}
    @{Row=9; I=6; J=6}
    @{Row=10; I=7; J=8}
    @{Row=11; I=8; J=8}
    @{Row=12; I=8; J=8}
    @{Row=13; I=10; J=10}
    @{Row=14; I=8; J=8}
    @{Row=15; I=10; J=10}
    @{Row=16; I=8; J=8}
    @{Row=17; I=8; J=8}
    @{Row=18; I=8; J=8}
    @{Row=19; I=8; J=8}
    @{Row=20; I=8; J=9}
    @{Row=21; I=8; J=8}
    @{Row=22; I=8; J=8}
    @{Row=23; I=8; J=8}
    @{Row=24; I=8; J=8}
    @{Row=25; I=7; J=8}
    @{Row=26; I=8; J=8}
    @{Row=27; I=8; J=8}
    @{Row=28; I=8; J=8}
    @{Row=29; I=8; J=9}
    @{Row=30; I=8; J=8}
    @{Row=31; I=8; J=8}
    @{Row=32; I=7; J=7}
    @{Row=33; I=11; J=11}
    @{Row=34; I=8; J=8}
    @{Row=35; I=7; J=8}
    @{Row=36; I=8; J=8}
    @{Row=37; I=7; J=8}
    @{Row=38; I=7; J=8}
    @{Row=39; I=8; J=8}
    @{Row=40; I=7; J=8}
    @{Row=41; I=7; J=7}
    @{Row=42; I=7; J=8}
    @{Row=43; I=7; J=8}
    @{Row=44; I=8; J=8}
    @{Row=45; I=9; J=9}
    @{Row=46; I=8; J=8}
    @{Row=47; I=8; J=8}
    @{Row=48; I=7; J=8}
    @{Row=49; I=8; J=8}
    @{Row=50; I=8; J=8}
    @{Row=51; I=8; J=8}
    @{Row=52; I=8; J=8}
    @{Row=53; I=7; J=7}
    @{Row=54; I=8; J=9}
    @{Row=55; I=8; J=8}
    @{Row=56; I=7; J=7}
    @{Row=57; I=4; J=4}
    @{Row=58; I=7; J=7}
)

foreach ($entry in $rowsData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 9).Value = $entry.I
    $ws.Cells.Item($r, 10).Value = $entry.J
}
